$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.079.37"
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = "3.544.46"
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'602.50"
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").Value = "'143.34"
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("D7").Value = "3.543.20"
$ws.Range("E7").Value = '  +1.09%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = "'0.489"
$ws.Range("E9").Value = '  +2.50%  '
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("D11").Value = "'7.76"
$ws.Range("E11").Value = '  -3.34%  '
$ws.Range("D12").Value = "'0.410"
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("D13").Value = "4.150.67"
$ws.Range("E13").Value = '  +1.19%  '
$ws.Range("D14").Value = "'0.0000204"
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").Value = "'29.87"
$ws.Range("E15").Value = '  -1.37%  '
$ws.Range("D16").Value = "3.536.95"
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = "66.152.96"
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").Value = "'11.31"
$ws.Range("E19").Value = '  +5.66%  '
$ws.Range("D20").Value = "'6.16"
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").Value = "'14.59"
$ws.Range("E21").Value = '  -1.60%  '
$ws.Range("D22").Value = "'428.22"
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D23").Value = "'0.605"
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("D24").Value = "'79.77"
$ws.Range("E24").Value = '  +2.26%  '
$ws.Range("D25").Value = "3.687.91"
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = "'0.0000115"
$ws.Range("E27").Value = '  -2.29%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = "'9.03"
$ws.Range("E29").Value = '  -2.65%  '
$ws.Range("D30").Value = "'7.77"
$ws.Range("E30").Value = '  -2.36%  '
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("D32").Value = "3.543.31"
$ws.Range("E32").Value = '  +1.39%  '
$ws.Range("D33").Value = "'25.28"
$ws.Range("E33").Value = '  +0.34%  '
$ws.Range("D34").Value = "'1.43"
$ws.Range("E34").Value = '  -2.09%  '
$ws.Range("E35").Value = '  -9.41%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = "'7.75"
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("D38").Value = "'1.70"
$ws.Range("E38").Value = '  -1.46%  '
$ws.Range("E39").Value = '  -1.87%  '
$ws.Range("D40").Value = "'174.29"
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("D41").Value = "'0.0841"
$ws.Range("E41").Value = '  -1.97%  '
$ws.Range("D42").Value = "'5.15"
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").Value = "'0.883"
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").Value = "'45.91"
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").Value = "'1.17"
$ws.Range("E47").Value = '  -3.07%  '
$ws.Range("D48").Value = "'24.54"
$ws.Range("E48").Value = '  -5.26%  '
$ws.Range("D49").Value = "'2.37"
$ws.Range("E49").Value = '  -2.32%  '
$ws.Range("D50").Value = "'7.07"
$ws.Range("E50").Value = '  -1.24%  '
$ws.Range("D51").Value = "'22.69"
$ws.Range("E51").Value = '  +0.35%  '
